$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing columns (headers "P" and "R", values "xxxx")
# that are no longer part of the model -- this also drops their now-unused
# shared strings and shrinks the sheet's dimension from AS7 to AQ7.
$ws.Range("AR1:AS7").EntireColumn.Delete()

# The "hVent" column (AO) used a placeholder value of 30; update it to the
# real coefficient of 0.03 for every data row.
$ws.Range("AO2:AO6").Value = 0.03

# Update the view: scroll right so column AE is the leftmost visible column
# and select AO6 (matches the saved window state in the edited workbook).
$excel.ActiveWindow.ScrollColumn = $ws.Range("AE1").Column
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AO6").Select()
